$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 849
$ws1.Range("F6").Value = 421
$ws1.Range("F7").Value = 612
$ws1.Range("F12").Value = 711
$ws1.Range("F14").Value = 1836
$ws1.Range("F16").Value = 3869
$ws1.Range("F17").Value = 377
$ws1.Range("F19").Value = 9
$ws1.Range("F20").Value = 60

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 479
$ws2.Range("F20").Value = 16

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5372
$ws3.Range("F4").Value = 289

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5372
$ws4.Range("F6").Value = 289
$ws4.Range("F12").Value = 479
$ws4.Range("F13").Value = 479
$ws4.Range("F14").Value = 849
$ws4.Range("F18").Value = 421
$ws4.Range("F19").Value = 612
$ws4.Range("F27").Value = 711
$ws4.Range("F30").Value = 1836
$ws4.Range("F32").Value = 3869
$ws4.Range("F34").Value = 377
$ws4.Range("F36").Value = 9
$ws4.Range("F37").Value = 60
$ws4.Range("F45").Value = 16
